# Applies the "ceklis berapa kriteria dari submission" edit:
#  - marks most of the checklist rows in column B with "V"
#  - adds a new "Fitur tambahan" block in columns AH/AI
#  - updates sheet view (top-left cell / selection) to match the new area

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B check marks ("V") on the existing checklist rows ---
$checkedRows = 5, 6, 7, 8, 11, 12, 13, 14, 16, 19, 21, 23
foreach ($r in $checkedRows) {
    $ws.Range("B$r").Value = "V"
}

# --- New "Fitur tambahan" (additional features) block ---
$ws.Range("AH2").Value = "Fitur tambahan"
$ws.Range("AI3").Value = "tambahkan kolom komentar"
$ws.Range("AI4").Value = "tambahkan tentang website ini"
$ws.Range("AI5").Value = "tambahkan list artikel"

# --- Column width for the new AH column, matching column O's width ---
$ws.Columns.Item(34).ColumnWidth = $ws.Columns.Item(15).ColumnWidth

# --- Sheet view: scroll / selection moved toward the newly added columns ---
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("AI6").Select()
